$d = $word.ActiveDocument

# -------------------------------------------------------------------
# "fix lab 4": the title run that reads "Bài thực hành số 4" should
# read "...5" instead. The number lives in its own run (separate from
# the "Bài thực hành số " run), so locate it with Find and then touch
# only that single-character run's Range -- this keeps the run's own
# rPr (font/lang) intact instead of merging the two runs together the
# way a Find-and-Replace across the run boundary would.
# -------------------------------------------------------------------
$titleScan = $d.Content
$titleScan.Find.Execute("Bài thực hành số ") | Out-Null
$numberPos = $titleScan.End
$numberRange = $d.Range($numberPos, $numberPos + 1)
$numberRange.Text = "5"

# -------------------------------------------------------------------
# "add lab 5": Word's "_GoBack" bookmark (last-edit marker) moves from
# where the previous editing session left it (end of the document,
# after the last paragraph) to right after the number we just edited
# -- exactly where Word drops it after you type/change text there.
# -------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Build the new (empty/collapsed) bookmark right after the "5". Adding
# a bookmark directly over a zero-length Range sitting on a paragraph
# boundary isn't reliable, so insert a throwaway character, bookmark
# the one-character range around it, then delete the character again
# -- the bookmark collapses to the empty range that is left behind.
$goBackPos = $numberRange.End
$goBackRange = $d.Range($goBackPos, $goBackPos)
$goBackRange.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $goBackRange)
$goBackRange.Text = ""

# -------------------------------------------------------------------
# Mark a handful of styles as "Quick Styles" (w:qFormat) -- Default
# Paragraph Font, Normal Table, and the WPS auto-generated TOC style.
# -------------------------------------------------------------------
$d.Styles("Default Paragraph Font").QuickStyle = $true
$d.Styles("Normal Table").QuickStyle = $true
$d.Styles("WPSOffice手动目录 1").QuickStyle = $true
